# -----------------------------------------------------------------------
# Commit: "Thu, Mar 19, 2020 11:06:19 AM"
#
# 1. The table on slide 16 gets a new PowerPoint built-in table style
#    (tableStyleId {33D33F59-...} -> {7372C059-...}).
# 2. The deck's theme (ppt/theme/theme1.xml, the one the slide master /
#    presentation point at) is switched from the custom "Integral"
#    colour scheme over to the standard Office theme colour scheme
#    (dk1/lt1 are already black/white in both, so only the other ten
#    slots actually change value).
# -----------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Re-style the financial table on slide 16 ------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{7372C059-4637-4A31-8576-EC7C662C14BC}")
    }
}

# --- 2. Swap the theme colour scheme from "Integral" to "Office" --------
function Pack-RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Index order exposed by ThemeColorScheme: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$officeThemeColors = @(
    (Pack-RGB 0x00 0x00 0x00),   # dk1       000000 (unchanged)
    (Pack-RGB 0xFF 0xFF 0xFF),   # lt1       FFFFFF (unchanged)
    (Pack-RGB 0x44 0x54 0x6A),   # dk2       44546A
    (Pack-RGB 0xE7 0xE6 0xE6),   # lt2       E7E6E6
    (Pack-RGB 0x5B 0x9B 0xD5),   # accent1   5B9BD5
    (Pack-RGB 0xED 0x7D 0x31),   # accent2   ED7D31
    (Pack-RGB 0xA5 0xA5 0xA5),   # accent3   A5A5A5
    (Pack-RGB 0xFF 0xC0 0x00),   # accent4   FFC000
    (Pack-RGB 0x44 0x72 0xC4),   # accent5   4472C4
    (Pack-RGB 0x70 0xAD 0x47),   # accent6   70AD47
    (Pack-RGB 0x05 0x63 0xC1),   # hlink     0563C1
    (Pack-RGB 0x95 0x4F 0x72)    # folHlink  954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeThemeColors[$i - 1]
}
